$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Friendly")

# Update the SE (speed-efficiency?) values in column F from 100 to 90
$ws.Range("F2").Value = 90
$ws.Range("F3").Value = 90
$ws.Range("F4").Value = 90

# Update the selection to span F2:F4 with F2 as the active cell
$ws.Range("F2:F4").Select()
